# Atualiza planilhas mensais (run 2025-09-18)
# Sheet "Precos": rename price column to "..._x", add a duplicate price
# column "..._y" (mirrors a pandas merge on "Nome do Produto" that leaves
# "_x"/"_y" suffixes), and re-sort the data rows by product name using a
# plain ordinal (codepoint) comparison - i.e. the same ordering Python's
# default string comparison / pandas sort_values would produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Precos")

# --- headers -----------------------------------------------------------
$ws.Range("B1").Value = "Preço_20250918_x"

$c1 = $ws.Range("C1")
$c1.Value = "Preço_20250918_y"
# Match the bold / centered / top-aligned / thin-bordered header style
# already used by A1 and B1.
$c1.Font.Bold = $true
$c1.HorizontalAlignment = -4108
$c1.VerticalAlignment = -4160
$c1.Borders.LineStyle = 1

# --- read the existing data rows ---------------------------------------
$dataRange = $ws.Range("A2:B36")
$vals = $dataRange.Value()
$n = $vals.GetLength(0)

$rows = @()
$i = 1
while ($i -le $n) {
    $rows += ,@($vals[$i,1], $vals[$i,2])
    $i = $i + 1
}

# --- sort rows by product name, ordinal/codepoint comparison -----------
# (.CompareTo() on [string] is an ordinal comparison here, matching
# Python's default string ordering - NOT Excel's usual locale/case
# insensitive Range.Sort order.)
$i = 1
while ($i -lt $rows.Length) {
    $j = $i
    while ($j -gt 0 -and ($rows[$j][0].CompareTo($rows[$j - 1][0]) -lt 0)) {
        $tmp = $rows[$j]
        $rows[$j] = $rows[$j - 1]
        $rows[$j - 1] = $tmp
        $j = $j - 1
    }
    $i = $i + 1
}

# --- write the sorted names/prices back, plus the duplicated "_y" col --
$k = 0
while ($k -lt $rows.Length) {
    $r = $k + 2
    $ws.Cells.Item($r, 1).Value = $rows[$k][0]
    $ws.Cells.Item($r, 2).Value = $rows[$k][1]
    $ws.Cells.Item($r, 3).Value = $rows[$k][1]
    $k = $k + 1
}
